$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.180.35'
$ws.Range('D2').Style = $origStyle

$ws.Range('E2').Value = '  +2.68%  '

$origStyle = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.199.76'
$ws.Range('D3').Style = $origStyle

$ws.Range('E3').Value = '  +1.19%  '

$origStyle = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = $origStyle

$ws.Range('E4').Value = '  -0.06%  '

$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '538.72'
$ws.Range('D5').Style = $origStyle

$ws.Range('E5').Value = '  +1.66%  '

$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.97'
$ws.Range('D6').Style = $origStyle

$ws.Range('E6').Value = '  +4.35%  '

$ws.Range('E7').Value = '  -0.04%  '

$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.530'
$ws.Range('D8').Style = $origStyle

$ws.Range('E8').Value = '  -2.21%  '

$ws.Range('E9').Value = '  +0.56%  '

$ws.Range('E10').Value = '  +1.15%  '

$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.433'
$ws.Range('D11').Style = $origStyle

$ws.Range('E11').Value = '  -1.38%  '

$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.748.89'
$ws.Range('D12').Style = $origStyle

$ws.Range('E12').Value = '  +1.11%  '

$ws.Range('E13').Value = '  -2.10%  '

$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.83'
$ws.Range('D14').Style = $origStyle

$ws.Range('E14').Value = '  +0.28%  '

$ws.Range('E15').Value = '  +0.76%  '

$origStyle = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '60.156.30'
$ws.Range('D16').Style = $origStyle

$ws.Range('E16').Value = '  +2.57%  '

$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.196.39'
$ws.Range('D17').Style = $origStyle

$ws.Range('E17').Value = '  +0.56%  '

$ws.Range('E18').Value = '  -0.06%  '

$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.30'
$ws.Range('D19').Style = $origStyle

$ws.Range('E19').Value = '  +2.32%  '

$ws.Range('E20').Value = '  +1.26%  '

$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '371.09'
$ws.Range('D21').Style = $origStyle

$ws.Range('E21').Value = '  -1.57%  '

$ws.Range('E22').Value = '  +0.13%  '

$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.524'
$ws.Range('D23').Style = $origStyle

$ws.Range('E23').Value = '  -1.46%  '

$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.63'
$ws.Range('D24').Style = $origStyle

$ws.Range('E24').Value = '  -0.13%  '

$ws.Range('E25').Value = '  +1.70%  '

$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.61'
$ws.Range('D26').Style = $origStyle

$ws.Range('E26').Value = '  +4.44%  '

$ws.Range('E27').Value = '  -0.04%  '

$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0879'
$ws.Range('D28').Style = $origStyle

$ws.Range('E28').Value = '  +1.22%  '

$ws.Range('E29').Value = '  +0.70%  '

$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.90'
$ws.Range('D30').Style = $origStyle

$ws.Range('E30').Value = '  +0.54%  '

$ws.Range('E31').Value = '  +1.13%  '

$ws.Range('E32').Value = '  +2.21%  '

$ws.Range('E34').Value = '  +3.01%  '

$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '159.07'
$ws.Range('D35').Style = $origStyle

$ws.Range('E35').Value = '  +1.35%  '

$ws.Range('E36').Value = '  +2.63%  '

$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '26.60'
$ws.Range('D37').Style = $origStyle

$ws.Range('E37').Value = '  +6.07%  '

$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.790.74'
$ws.Range('D38').Style = $origStyle

$ws.Range('E38').Value = '  +4.79%  '

$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0314'
$ws.Range('D39').Style = $origStyle

$ws.Range('E39').Value = '  +8.59%  '

$ws.Range('E40').Value = '  +2.05%  '

$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.70'
$ws.Range('D41').Style = $origStyle

$ws.Range('E41').Value = '  +0.54%  '

$ws.Range('E42').Value = '  -1.33%  '

$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '39.89'
$ws.Range('D43').Style = $origStyle

$ws.Range('E43').Value = '  +2.06%  '

$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.719'
$ws.Range('D44').Style = $origStyle

$ws.Range('E44').Value = '  -0.45%  '

$ws.Range('E45').Value = '  +0.85%  '

$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.237.09'
$ws.Range('D46').Style = $origStyle

$ws.Range('E46').Value = '  +1.05%  '

$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.984'
$ws.Range('D47').Style = $origStyle

$ws.Range('E47').Value = '  +0.18%  '

$ws.Range('E48').Value = '  -0.88%  '

$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '20.67'
$ws.Range('D49').Style = $origStyle

$ws.Range('E49').Value = '  +2.81%  '

$ws.Range('E50').Value = '  +5.78%  '

$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.999'
$ws.Range('D51').Style = $origStyle

$ws.Range('E51').Value = '  -0.05%  '
